$d = $word.ActiveDocument

# Update the date line (first paragraph, outside the table)
$d.Paragraphs.Item(1).Range.Text = "2024-09-23 Monday"

# Update each table cell value in row-major order (5 columns per row)
$t = $d.Tables.Item(1)
$values = @(
    "60-46=14",
    "26+13=39",
    "66-42=24",
    "74+5=79",
    "21+4=25",
    "23-7=16",
    "1+46=47",
    "71-22=49",
    "85-1=84",
    "76-23=53",
    "26+60=86",
    "80-68=12",
    "64-42=22",
    "89-86=3",
    "94-43=51",
    "48+50=98",
    "48+4=52",
    "50+24=74",
    "45-42=3",
    "79-35=44",
    "98-6=92",
    "83+14=97",
    "11+8=19",
    "37-12=25",
    "79+4=83",
    "84-56=28",
    "43+32=75",
    "77-10=67",
    "73-46=27",
    "88-79=9",
    "54+38=92",
    "62-5=57",
    "67-17=50",
    "98-90=8",
    "76-21=55",
    "34+40=74",
    "74-21=53",
    "36+30=66",
    "82-55=27",
    "96-96=0",
    "45-39=6",
    "32-30=2",
    "5+10=15",
    "14-8=6",
    "62+26=88",
    "17+60=77",
    "32+16=48",
    "41+1=42",
    "68+27=95",
    "41-33=8",
    "31+11=42",
    "31+25=56",
    "36-12=24",
    "27-15=12",
    "9+51=60",
    "42+36=78",
    "88-61=27",
    "73+21=94",
    "16+22=38",
    "16+74=90",
    "52-49=3",
    "25-21=4",
    "57-9=48",
    "32+35=67",
    "92-46=46",
    "80-33=47",
    "15+54=69",
    "61-13=48",
    "11+63=74",
    "88-13=75",
    "2+0=2",
    "35+5=40",
    "47+45=92",
    "4+3=7",
    "12+71=83",
    "96-73=23",
    "60+7=67",
    "1+0=1",
    "90-14=76",
    "62+6=68",
    "35-13=22",
    "2+48=50",
    "94-10=84",
    "41+21=62",
    "58+30=88",
    "44-41=3",
    "18+30=48",
    "26+60=86",
    "97-58=39",
    "17+50=67",
    "34-25=9",
    "71+13=84",
    "97-18=79",
    "15+44=59",
    "69-2=67",
    "58-44=14",
    "32+20=52",
    "58+26=84",
    "78-8=70",
    "95-52=43"
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $t.Cell($r, $c).Range.Text = $values[$idx]
        $idx++
    }
}

Write-Host "Done. Updated" $idx "cells."